$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text
# (matching the workbook convention of storing Price/Volume figures as strings).
# Set NumberFormat to Text ("@") first so Excel does not coerce them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D5").Value = "552.05"
$ws.Range("D6").Value = "137.16"
$ws.Range("D10").Value = "5.67"
$ws.Range("D11").Value = "0.147"
$ws.Range("D13").Value = "25.39"
$ws.Range("D18").Value = "11.30"
$ws.Range("D20").Value = "329.21"
$ws.Range("D21").Value = "6.65"
$ws.Range("D23").Value = "65.63"
$ws.Range("D24").Value = "0.178"
$ws.Range("D25").Value = "8.61"
$ws.Range("D30").Value = "168.91"
$ws.Range("D31").Value = "6.05"
$ws.Range("D32").Value = "18.59"
$ws.Range("D33").Value = "1.02"
$ws.Range("D37").Value = "4.18"
$ws.Range("D39").Value = "321.06"
$ws.Range("D41").Value = "3.66"
$ws.Range("D42").Value = "139.91"
$ws.Range("D45").Value = "0.0514"
$ws.Range("D46").Value = "0.576"
$ws.Range("D48").Value = "0.386"

# Remaining cells: values that will naturally stay text (contain non-numeric characters)
$ws.Range("D2").Value = "60.017.51"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.417.29"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +4.08%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "2.847.50"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "59.938.64"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.418.78"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("E51").Value = "  -1.09%  "
